$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# New "Status" text for the ee8c9196... row (row 7) on both language sheets,
# and the matching summary cells on the Overview sheet (same underlying text).
$overview.Range("E7").Value = "Handback transform failed"
$overview.Range("F7").Value = "Handback transform failed"
$zhcn.Range("C7").Value = "Handback transform failed"
$dede.Range("C7").Value = "Handback transform failed"

# New "Error Detail" text (column P) for the same row on both language sheets.
$zhcn.Range("P7").Value = "Handback file name: uhonibjd.spy is different with handoff file name: ee8c9196-0a69-4ee2-8c6d-7d66b6c6b10e.20feb19b3ab7bcfc39302a5158b070200ffeb913.zh-cn."
$dede.Range("P7").Value = "Handback file name: uhonibjd.spy is different with handoff file name: ee8c9196-0a69-4ee2-8c6d-7d66b6c6b10e.20feb19b3ab7bcfc39302a5158b070200ffeb913.de-de."

# Widen the "Error Detail" column to fit the new, longer text (stored XML
# width = ColumnWidth + 5/6, so back the padding out to land on exactly 40).
$zhcn.Columns.Item(16).ColumnWidth = 39.1666666666667
$dede.Columns.Item(16).ColumnWidth = 39.1666666666667
